# Auto-generated edit script applying the scheduled-runner price/profit refresh
# to the Exodus_Profits workbook. Each sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# gets a batch of cell writes; cells that must disappear entirely (blank cache slots)
# are cleared with ClearContents so they are dropped from the saved XML, matching the diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 56.25
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H18").Value = 749.8333
$ws.Range("I18").Value = 749.8333
$ws.Range("K18").Value = 749.8333
$ws.Range("M18").Value = -465.8333
$ws.Range("H33").Value = 204.2
$ws.Range("J33").Value = 850
$ws.Range("L33").Value = 850
$ws.Range("N33").Value = -1308
$ws.Range("H86").Value = 3691.4119
$ws.Range("I86").Value = 3260.8333
$ws.Range("J86").Value = 3926.2727
$ws.Range("K86").Value = 3260.8333
$ws.Range("L86").Value = 3926.2727
$ws.Range("M86").Value = -2137.8333
$ws.Range("N86").Value = -6172.2727
$ws.Range("H89").Value = 3691.4119
$ws.Range("I89").Value = 3260.8333
$ws.Range("J89").Value = 3926.2727
$ws.Range("K89").Value = 16304.1665
$ws.Range("L89").Value = 19631.3635
$ws.Range("M89").Value = -10688.1665
$ws.Range("N89").Value = -30863.3635
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H137").Value = 662261.1
$ws.Range("J137").Value = 1119029.9
$ws.Range("L137").Value = 3357089.7
$ws.Range("N137").Value = -3362189.7
$ws.Range("H138").Value = 1695.7727
$ws.Range("I138").Value = 1320.8
$ws.Range("J138").Value = 2499.2856
$ws.Range("K138").Value = 3962.4
$ws.Range("L138").Value = 7497.8568
$ws.Range("M138").Value = 1177.6
$ws.Range("N138").Value = -17777.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 88998
$ws.Range("J7").Value = 88998
$ws.Range("L7").Value = 88998
$ws.Range("N7").Value = -89226
$ws.Range("H53").Value = 1039
$ws.Range("I53").Value = 1039
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 1039
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -357
$ws.Range("N53").ClearContents()
$ws.Range("H74").Value = 35683.8
$ws.Range("I74").Value = 51471.55
$ws.Range("K74").Value = 51471.55
$ws.Range("M74").Value = -50597.55
$ws.Range("H77").Value = 35683.8
$ws.Range("I77").Value = 51471.55
$ws.Range("K77").Value = 257357.75
$ws.Range("M77").Value = -252989.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 106963.77
$ws.Range("J20").Value = 3466.3333
$ws.Range("L20").Value = 3466.3333
$ws.Range("N20").Value = -3960.3333
$ws.Range("H86").Value = 3230.05
$ws.Range("I86").Value = 3464.6428
$ws.Range("J86").Value = 2682.6667
$ws.Range("K86").Value = 3464.6428
$ws.Range("L86").Value = 2682.6667
$ws.Range("M86").Value = -2341.6428
$ws.Range("N86").Value = -4928.6667
$ws.Range("H89").Value = 3230.05
$ws.Range("I89").Value = 3464.6428
$ws.Range("J89").Value = 2682.6667
$ws.Range("K89").Value = 17323.214
$ws.Range("L89").Value = 13413.3335
$ws.Range("M89").Value = -11707.214
$ws.Range("N89").Value = -24645.3335
$ws.Range("H107").Value = 8335646.5
$ws.Range("I107").Value = 11113233
$ws.Range("K107").Value = 11113233
$ws.Range("M107").Value = -11111313

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1446.4
$ws.Range("I16").Value = 1065.8462
$ws.Range("K16").Value = 1065.8462
$ws.Range("M16").Value = -778.8462
$ws.Range("H69").Value = 45999.332
$ws.Range("I69").Value = 45999.332
$ws.Range("K69").Value = 45999.332
$ws.Range("M69").Value = -45250.332
$ws.Range("H72").Value = 45999.332
$ws.Range("I72").Value = 45999.332
$ws.Range("K72").Value = 137997.996
$ws.Range("M72").Value = -134253.996
$ws.Range("H99").Value = 14238416
$ws.Range("I99").Value = 15875380
$ws.Range("K99").Value = 15875380
$ws.Range("M99").Value = -15873882
$ws.Range("H113").Value = 1446.4
$ws.Range("I113").Value = 1065.8462
$ws.Range("K113").Value = 1065.8462
$ws.Range("M113").Value = 1104.1538
$ws.Range("H126").Value = 14238416
$ws.Range("I126").Value = 15875380
$ws.Range("K126").Value = 47626140
$ws.Range("M126").Value = -47623670

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 85.818184
$ws.Range("I2").Value = 158.66667
$ws.Range("J2").Value = 35.384617
$ws.Range("K2").Value = 952.0000200000001
$ws.Range("L2").Value = 212.307702
$ws.Range("M2").Value = -839.0000200000001
$ws.Range("N2").Value = -438.307702
$ws.Range("H4").Value = 31666838
$ws.Range("I4").Value = 33333334
$ws.Range("J4").Value = 30000342
$ws.Range("K4").Value = 100000002
$ws.Range("L4").Value = 90001026
$ws.Range("M4").Value = -99999890
$ws.Range("N4").Value = -90001250
$ws.Range("H7").Value = 3578.6667
$ws.Range("J7").Value = 6940.3335
$ws.Range("L7").Value = 20821.0005
$ws.Range("N7").Value = -21045.0005
$ws.Range("H9").Value = 15100
$ws.Range("I9").Value = 15100
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 45300
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -45076
$ws.Range("N9").ClearContents()
$ws.Range("H11").Value = 2972.25
$ws.Range("I11").Value = 5669.75
$ws.Range("J11").Value = 274.75
$ws.Range("K11").Value = 17009.25
$ws.Range("L11").Value = 824.25
$ws.Range("M11").Value = -16869.25
$ws.Range("N11").Value = -1104.25
$ws.Range("H23").Value = 72718.14
$ws.Range("I23").Value = 929.5
$ws.Range("J23").Value = 101433.6
$ws.Range("K23").Value = 2788.5
$ws.Range("L23").Value = 304300.8
$ws.Range("M23").Value = -2553.5
$ws.Range("N23").Value = -304770.8
$ws.Range("H25").Value = 25000218
$ws.Range("I25").Value = 275
$ws.Range("J25").Value = 33333534
$ws.Range("K25").Value = 825
$ws.Range("L25").Value = 100000602
$ws.Range("M25").Value = -656
$ws.Range("N25").Value = -100000940
$ws.Range("H30").Value = 25000218
$ws.Range("I30").Value = 275
$ws.Range("J30").Value = 33333534
$ws.Range("K30").Value = 825
$ws.Range("L30").Value = 100000602
$ws.Range("M30").Value = -723
$ws.Range("N30").Value = -100000806
$ws.Range("H110").Value = 6885.4
$ws.Range("J110").Value = 7333.3335
$ws.Range("L110").Value = 22000.0005
$ws.Range("N110").Value = -30180.0005
$ws.Range("H129").Value = 47619536
$ws.Range("I129").Value = 568.8333
$ws.Range("K129").Value = 1706.4999
$ws.Range("M129").Value = 3293.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 289.8421
$ws.Range("I2").Value = 300.08334
$ws.Range("J2").Value = 272.2857
$ws.Range("K2").Value = 300.08334
$ws.Range("L2").Value = 272.2857
$ws.Range("M2").Value = -187.08334
$ws.Range("N2").Value = -498.2857
$ws.Range("J70").Value = 11000
$ws.Range("L70").Value = 11000
$ws.Range("N70").Value = -11540
$ws.Range("J73").Value = 11000
$ws.Range("L73").Value = 11000
$ws.Range("N73").Value = -12872
$ws.Range("H92").Value = 4000
$ws.Range("J92").Value = 4000
$ws.Range("L92").Value = 4000
$ws.Range("N92").Value = -7744
$ws.Range("H113").Value = 4923508
$ws.Range("I113").Value = 372740.66
$ws.Range("J113").Value = 8336583
$ws.Range("K113").Value = 372740.66
$ws.Range("L113").Value = 8336583
$ws.Range("M113").Value = -370570.66
$ws.Range("N113").Value = -8340923
$ws.Range("H132").Value = 4096.5654
$ws.Range("I132").Value = 3222.158
$ws.Range("K132").Value = 9666.474
$ws.Range("M132").Value = -7136.474

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1063.0834
$ws.Range("I22").Value = 1201
$ws.Range("J22").Value = 833.2222
$ws.Range("K22").Value = 1201
$ws.Range("L22").Value = 833.2222
$ws.Range("M22").Value = -906
$ws.Range("N22").Value = -1423.2222
$ws.Range("H27").Value = 1063.0834
$ws.Range("I27").Value = 1201
$ws.Range("J27").Value = 833.2222
$ws.Range("K27").Value = 1201
$ws.Range("L27").Value = 833.2222
$ws.Range("M27").Value = -1094
$ws.Range("N27").Value = -1047.2222
$ws.Range("H40").Value = 9154532
$ws.Range("I40").Value = 4110.778
$ws.Range("K40").Value = 4110.778
$ws.Range("M40").Value = -3974.778
$ws.Range("H45").Value = 13298.143
$ws.Range("I45").Value = 18347
$ws.Range("J45").Value = 9511.5
$ws.Range("K45").Value = 18347
$ws.Range("L45").Value = 9511.5
$ws.Range("M45").Value = -17940
$ws.Range("N45").Value = -10325.5
$ws.Range("H46").Value = 1625.1538
$ws.Range("J46").Value = 1800
$ws.Range("L46").Value = 1800
$ws.Range("N46").Value = -2176
$ws.Range("H55").Value = 4762535
$ws.Range("I55").Value = 453.35294
$ws.Range("J55").Value = 8000750.5
$ws.Range("K55").Value = 453.35294
$ws.Range("L55").Value = 8000750.5
$ws.Range("M55").Value = -280.35294
$ws.Range("N55").Value = -8001096.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9709
$ws.Range("I41").Value = 6750
$ws.Range("J41").Value = 11188.5
$ws.Range("K41").Value = 6750
$ws.Range("L41").Value = 11188.5
$ws.Range("M41").Value = -6360
$ws.Range("N41").Value = -11968.5
$ws.Range("H126").Value = 81721.78999999999
$ws.Range("I126").Value = 102259.55
$ws.Range("J126").Value = 6416.6665
$ws.Range("K126").Value = 306778.65
$ws.Range("L126").Value = 19249.9995
$ws.Range("M126").Value = -304308.65
$ws.Range("N126").Value = -24189.9995

Write-Output "Applied scheduled-runner profit refresh across all sheets."
